# Update the "想去人数" (column F) values on the "展览" and "全部类型"
# worksheets to reflect the latest scrape counts.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 198
    3  = 768
    6  = 4543
    8  = 366
    9  = 1310
    10 = 540
    12 = 901
    14 = 505
    15 = 55
    16 = 238
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
